$wb = $excel.ActiveWorkbook

# --- "Week 22" sheet: fill in attendance hours (Tue-Fri) for every person ---
$week22 = $wb.Worksheets.Item("Week 22")

$week22.Range("B3:E3").Value = 4   # Tue
$week22.Range("B4:E4").Value = 2   # Wed
$week22.Range("B5:E5").Value = 4   # Thur
$week22.Range("B6:E6").Value = 8   # Fri

# Remember the author's last clicked cell on this sheet
$week22.Range("E6").Select()

# --- "Week 23" sheet: no data entered yet, only the remembered selection moved ---
$week23 = $wb.Worksheets.Item("Week 23")
$week23.Range("E32").Select()

# Re-select on the originally active sheet ("Percentages") so it stays the active tab
$pct = $wb.Worksheets.Item("Percentages")
$pct.Range("C6").Select()
